$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 40; this shifts the existing rows 40-70
# down to 41-71 (carrying their values/styles along), matching the diff's
# "row N becomes row N+1" pattern for rows 40..70.
$ws.Rows.Item(40).Insert()

# Populate the newly inserted row 40 with the new weekly data point.
$ws.Range("A40").Value = 9
$ws.Range("B40").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C40").Value = "Metropolitana"
$ws.Range("D40").Value = 44827
$ws.Range("E40").Value = 13
$ws.Range("F40").Value = 100112035
$ws.Range("G40").Value = "Bruselas (repollito)"
$ws.Range("H40").Value = "Sin especificar"
$ws.Range("I40").Value = "Primera"
$ws.Range("J40").Value = 45
$ws.Range("K40").Value = 19000
$ws.Range("L40").Value = 20000
$ws.Range("M40").Value = 19556
$ws.Range("N40").Value = "$/malla 15 kilos"
$ws.Range("O40").Value = "Hijuelas"
$ws.Range("P40").Value = 1304
$ws.Range("Q40").Value = 15
$ws.Range("R40").Value = "Hortaliza"
